$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Cells.Item(2, 11).Value = 5970
$ws.Cells.Item(3, 11).Value = 6143
$ws.Cells.Item(4, 11).Value = 1289
$ws.Cells.Item(5, 11).Value = 438
$ws.Cells.Item(6, 11).Value = 6765
$ws.Cells.Item(7, 11).Value = 20605

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Cells.Item(2, 11).Value = 181
$ws.Cells.Item(6, 11).Value = 144
$ws.Cells.Item(7, 11).Value = 598
$ws.Cells.Item(8, 11).Value = 1359
$ws.Cells.Item(10, 11).Value = 116
$ws.Cells.Item(15, 11).Value = 212
$ws.Cells.Item(18, 11).Value = 136
$ws.Cells.Item(19, 11).Value = 598
$ws.Cells.Item(22, 11).Value = 61
$ws.Cells.Item(23, 11).Value = 211
$ws.Cells.Item(29, 11).Value = 1125
$ws.Cells.Item(33, 11).Value = 892
$ws.Cells.Item(34, 11).Value = 117
$ws.Cells.Item(37, 11).Value = 698
$ws.Cells.Item(39, 11).Value = 25
$ws.Cells.Item(41, 11).Value = 143
$ws.Cells.Item(42, 11).Value = 767
$ws.Cells.Item(49, 11).Value = 113
$ws.Cells.Item(50, 11).Value = 99
$ws.Cells.Item(53, 11).Value = 264
$ws.Cells.Item(54, 11).Value = 401
$ws.Cells.Item(57, 11).Value = 78
$ws.Cells.Item(59, 11).Value = 35
$ws.Cells.Item(63, 11).Value = 64
$ws.Cells.Item(65, 11).Value = 485
$ws.Cells.Item(67, 11).Value = 805
$ws.Cells.Item(70, 11).Value = 35
$ws.Cells.Item(72, 11).Value = 100
$ws.Cells.Item(78, 11).Value = 235
$ws.Cells.Item(79, 11).Value = 515
$ws.Cells.Item(83, 11).Value = 455
$ws.Cells.Item(85, 11).Value = 956
$ws.Cells.Item(87, 11).Value = 35
$ws.Cells.Item(89, 11).Value = 301
$ws.Cells.Item(90, 11).Value = 190
$ws.Cells.Item(92, 11).Value = 79
$ws.Cells.Item(94, 11).Value = 276
$ws.Cells.Item(95, 11).Value = 348
$ws.Cells.Item(96, 11).Value = 214
$ws.Cells.Item(98, 11).Value = 100
$ws.Cells.Item(99, 11).Value = 338
$ws.Cells.Item(101, 11).Value = 20605

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Cells.Item(6, 11).Value = 93
$ws.Cells.Item(7, 11).Value = 214

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Cells.Item(3, 11).Value = 195
$ws.Cells.Item(6, 11).Value = 161
$ws.Cells.Item(7, 11).Value = 598

$ws = $wb.Worksheets.Item('Uptown')
$ws.Cells.Item(4, 11).Value = 36
$ws.Cells.Item(7, 11).Value = 301

$ws = $wb.Worksheets.Item('South Shore')
$ws.Cells.Item(3, 11).Value = 330
$ws.Cells.Item(4, 11).Value = 53
$ws.Cells.Item(7, 11).Value = 956

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Cells.Item(2, 11).Value = 67
$ws.Cells.Item(6, 11).Value = 114
$ws.Cells.Item(7, 11).Value = 264

$ws = $wb.Worksheets.Item('Austin')
$ws.Cells.Item(2, 11).Value = 374
$ws.Cells.Item(3, 11).Value = 416
$ws.Cells.Item(6, 11).Value = 456
$ws.Cells.Item(7, 11).Value = 1359

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Cells.Item(6, 11).Value = 105
$ws.Cells.Item(7, 11).Value = 455

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Cells.Item(3, 11).Value = 329
$ws.Cells.Item(4, 11).Value = 44
$ws.Cells.Item(5, 11).Value = 24
$ws.Cells.Item(6, 11).Value = 261
$ws.Cells.Item(7, 11).Value = 892

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Cells.Item(2, 11).Value = 119
$ws.Cells.Item(3, 11).Value = 123
$ws.Cells.Item(7, 11).Value = 348

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Cells.Item(3, 11).Value = 232
$ws.Cells.Item(4, 11).Value = 34
$ws.Cells.Item(7, 11).Value = 698

$ws = $wb.Worksheets.Item('New City')
$ws.Cells.Item(2, 11).Value = 160
$ws.Cells.Item(3, 11).Value = 119
$ws.Cells.Item(4, 11).Value = 18
$ws.Cells.Item(6, 11).Value = 176
$ws.Cells.Item(7, 11).Value = 485

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Cells.Item(2, 11).Value = 87
$ws.Cells.Item(7, 11).Value = 338

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Cells.Item(3, 11).Value = 287
$ws.Cells.Item(6, 11).Value = 230
$ws.Cells.Item(7, 11).Value = 805

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Cells.Item(6, 11).Value = 55
$ws.Cells.Item(7, 11).Value = 113

$ws = $wb.Worksheets.Item('Loop')
$ws.Cells.Item(3, 11).Value = 97
$ws.Cells.Item(6, 11).Value = 217
$ws.Cells.Item(7, 11).Value = 401

$ws = $wb.Worksheets.Item('Englewood')
$ws.Cells.Item(2, 11).Value = 322
$ws.Cells.Item(3, 11).Value = 403
$ws.Cells.Item(4, 11).Value = 55
$ws.Cells.Item(7, 11).Value = 1125

$ws = $wb.Worksheets.Item('Chatham')
$ws.Cells.Item(2, 11).Value = 177
$ws.Cells.Item(4, 11).Value = 29
$ws.Cells.Item(6, 11).Value = 192
$ws.Cells.Item(7, 11).Value = 598

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Cells.Item(4, 11).Value = 9
$ws.Cells.Item(7, 11).Value = 144

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Cells.Item(6, 11).Value = 57
$ws.Cells.Item(7, 11).Value = 143

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Cells.Item(2, 11).Value = 210
$ws.Cells.Item(5, 11).Value = 9
$ws.Cells.Item(6, 11).Value = 285
$ws.Cells.Item(7, 11).Value = 767

$ws = $wb.Worksheets.Item('Avondale')
$ws.Cells.Item(2, 11).Value = 35
$ws.Cells.Item(6, 11).Value = 53
$ws.Cells.Item(7, 11).Value = 116

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Cells.Item(2, 11).Value = 70
$ws.Cells.Item(7, 11).Value = 235

$ws = $wb.Worksheets.Item('Douglas')
$ws.Cells.Item(5, 11).Value = 9
$ws.Cells.Item(7, 11).Value = 211

$ws = $wb.Worksheets.Item('Roseland')
$ws.Cells.Item(2, 11).Value = 173
$ws.Cells.Item(3, 11).Value = 168
$ws.Cells.Item(6, 11).Value = 127
$ws.Cells.Item(7, 11).Value = 515

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Cells.Item(2, 11).Value = 162
$ws.Cells.Item(6, 11).Value = 138

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Cells.Item(3, 11).Value = 44
$ws.Cells.Item(6, 11).Value = 36
$ws.Cells.Item(7, 11).Value = 136

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Cells.Item(2, 11).Value = 46
$ws.Cells.Item(7, 11).Value = 117

$ws = $wb.Worksheets.Item('West Loop')
$ws.Cells.Item(2, 11).Value = 74
$ws.Cells.Item(6, 11).Value = 123
$ws.Cells.Item(7, 11).Value = 276

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Cells.Item(2, 11).Value = 75
$ws.Cells.Item(3, 11).Value = 53
$ws.Cells.Item(4, 11).Value = 13
$ws.Cells.Item(7, 11).Value = 212

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Cells.Item(2, 11).Value = 19
$ws.Cells.Item(7, 11).Value = 100

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Cells.Item(4, 11).Value = 12
$ws.Cells.Item(7, 11).Value = 99

$ws = $wb.Worksheets.Item('Greektown')
$ws.Cells.Item(5, 11).Value = 15
$ws.Cells.Item(6, 11).Value = 25

$ws = $wb.Worksheets.Item('Montclare')
$ws.Cells.Item(2, 11).Value = 10
$ws.Cells.Item(7, 11).Value = 35

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Cells.Item(2, 11).Value = 54
$ws.Cells.Item(7, 11).Value = 181

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Cells.Item(6, 11).Value = 38
$ws.Cells.Item(7, 11).Value = 79

$ws = $wb.Worksheets.Item('O''Hare')
$ws.Cells.Item(4, 11).Value = 6
$ws.Cells.Item(7, 11).Value = 35

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Cells.Item(2, 11).Value = 70
$ws.Cells.Item(7, 11).Value = 190

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Cells.Item(6, 11).Value = 35
$ws.Cells.Item(7, 11).Value = 78

$ws = $wb.Worksheets.Item('Clearing')
$ws.Cells.Item(2, 11).Value = 28
$ws.Cells.Item(7, 11).Value = 61

$ws = $wb.Worksheets.Item('Old Town')
$ws.Cells.Item(2, 11).Value = 18
$ws.Cells.Item(7, 11).Value = 100

$ws = $wb.Worksheets.Item('Ukrainian Village')
$ws.Cells.Item(2, 11).Value = 4
$ws.Cells.Item(7, 11).Value = 35
